# Updates cryptos list prices / 1h volume percentages (GitHub Actions refresh).
# Column layout: A=index, B=Coin, C=Link, D=Price, E=Volume(1h)
# Price values that look like plain decimal numbers are written with a
# leading apostrophe so Excel keeps them as text (matching the source
# workbook, which stores every Price/Volume cell as a string).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bitcoin
$ws.Range("D2").Value = "29.010.97"
$ws.Range("E2").Value = "  +0.29%  "

# Ethereum
$ws.Range("D3").Value = "1.921.35"
$ws.Range("E3").Value = "  +1.54%  "

# TetherUSD
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.15%  "

# BNB
$ws.Range("D5").Value = "'325.26"
$ws.Range("E5").Value = "  +0.47%  "

# USDC
$ws.Range("E6").Value = "  +0.21%  "

# XRP
$ws.Range("D7").Value = "'0.4595"
$ws.Range("E7").Value = "  +0.06%  "

# Cardano
$ws.Range("D8").Value = "'0.3816"
$ws.Range("E8").Value = "  -0.07%  "

# Dogecoin
$ws.Range("D9").Value = "'0.07748"
$ws.Range("E9").Value = "  +0.22%  "

# Polygon
$ws.Range("D10").Value = "'0.9775"
$ws.Range("E10").Value = "  +1.06%  "

# Solana
$ws.Range("D11").Value = "'22.76"
$ws.Range("E11").Value = "  +3.08%  "

# WrappedEther
$ws.Range("D12").Value = "1.924.24"
$ws.Range("E12").Value = "  +2.15%  "

# Polkadot
$ws.Range("D13").Value = "'5.697"
$ws.Range("E13").Value = "  +0.17%  "

# Chainlink
$ws.Range("D14").Value = "'6.966"
$ws.Range("E14").Value = "  +0.11%  "

# TRON
$ws.Range("D15").Value = "'0.07002"
$ws.Range("E15").Value = "  -0.36%  "

# BinanceUSD
$ws.Range("E16").Value = "  +0.25%  "

# Litecoin
$ws.Range("D17").Value = "'84.24"
$ws.Range("E17").Value = "  +0.90%  "

# ShibaInu
$ws.Range("D18").Value = "'0.000009521"
$ws.Range("E18").Value = "  -0.18%  "

# Avalanche
$ws.Range("D19").Value = "'16.69"
$ws.Range("E19").Value = "  +0.15%  "

# Dai
$ws.Range("E20").Value = "  +0.16%  "

# WrappedBTC
$ws.Range("D21").Value = "29.016.46"
$ws.Range("E21").Value = "  +0.50%  "

# Uniswap
$ws.Range("D22").Value = "'5.355"
$ws.Range("E22").Value = "  +0.82%  "

# Cosmos
$ws.Range("D23").Value = "'11.02"
$ws.Range("E23").Value = "  +0.95%  "

# WrappedliquidstakedEther2.0
$ws.Range("D24").Value = "2.155.28"
$ws.Range("E24").Value = "  +1.86%  "

# Toncoin
$ws.Range("D25").Value = "'2.073"
$ws.Range("E25").Value = "  +0.17%  "

# Monero
$ws.Range("D26").Value = "'157.27"
$ws.Range("E26").Value = "  +0.70%  "

# EthereumClassic
$ws.Range("D27").Value = "'19.06"
$ws.Range("E27").Value = "  +0.04%  "

# InternetComputer(DFINITY)
$ws.Range("E28").Value = "  +0.74%  "

# BitcoinCash
$ws.Range("D29").Value = "'117.60"
$ws.Range("E29").Value = "  +0.13%  "

# LidoDAOToken
$ws.Range("D30").Value = "'1.836"
$ws.Range("E30").Value = "  +1.02%  "

# Stellar
$ws.Range("D31").Value = "'0.09334"
$ws.Range("E31").Value = "  +0.75%  "

# ImmutableX
$ws.Range("D32").Value = "'0.8585"
$ws.Range("E32").Value = "  +0.68%  "

# Filecoin
$ws.Range("D33").Value = "'5.096"
$ws.Range("E33").Value = "  +0.28%  "

# ARBITRUM
$ws.Range("D34").Value = "'1.238"
$ws.Range("E34").Value = "  +0.01%  "

# HuobiToken
$ws.Range("D35").Value = "'3.016"
$ws.Range("E35").Value = "  +0.15%  "

# Row 36/37 swap places: TrustWalletToken <-> Hedera
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "'0.05682"
$ws.Range("E36").Value = "  -0.17%  "

$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "'1.158"
$ws.Range("E37").Value = "  +0.85%  "

# Frax
$ws.Range("D38").Value = "'1.002"
$ws.Range("E38").Value = "  +0.10%  "

# VeChain
$ws.Range("D39").Value = "'0.02047"
$ws.Range("E39").Value = "  +0.35%  "

# MXToken
$ws.Range("D40").Value = "'3.104"
$ws.Range("E40").Value = "  +15.08%  "

# FraxShare
$ws.Range("D41").Value = "'7.419"
$ws.Range("E41").Value = "  -0.04%  "

# TheSandbox
$ws.Range("D42").Value = "'0.5497"
$ws.Range("E42").Value = "  -0.05%  "

# Algorand
$ws.Range("D43").Value = "'0.1753"
$ws.Range("E43").Value = "  +0.17%  "

# Aptos
$ws.Range("D44").Value = "'9.359"
$ws.Range("E44").Value = "  +1.52%  "

# PEPE
$ws.Range("D45").Value = "'0.000002836"
$ws.Range("E45").Value = "  -1.52%  "

# RenderToken
$ws.Range("D46").Value = "'2.182"
$ws.Range("E46").Value = "  +5.08%  "

# Decentraland
$ws.Range("D47").Value = "'0.5196"
$ws.Range("E47").Value = "  +0.30%  "

# EnergySwap
$ws.Range("D48").Value = "'11.27"
$ws.Range("E48").Value = "  -0.18%  "

# Cronos
$ws.Range("D49").Value = "'0.06912"
$ws.Range("E49").Value = "  +1.55%  "

# Quant
$ws.Range("D50").Value = "'110.37"
$ws.Range("E50").Value = "  -1.06%  "

# NEARProtocol
$ws.Range("D51").Value = "'1.760"
$ws.Range("E51").Value = "  -1.24%  "
